# Update Excel workbook: Spain provincias COVID data
# - Refreshes the "last updated" timestamp in A1
# - Updates Casos totales / Casos activos / Recuperados / Muertes (cols B-E)
#   for the provinces whose figures changed in this update

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." banner
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 12:46"

# Row => @(Casos totales, Casos activos, Recuperados, Muertes)
$updates = @{
    7  = @(1501, 466, 1457, 44)
    8  = @(1497, 23, 1413, 61)
    9  = @(1207, 466, 1134, 73)
    10  = @(1197, 23, 1141, 33)
    12  = @(907, 22, 845, 40)
    13  = @(857, 12, 774, 71)
    14  = @(779, 35, 719, 25)
    15  = @(734, 25, 715, 19)
    16  = @(665, 72, 639, 26)
    17  = @(602, 25, 598, 4)
    18  = @(563, 466, 547, 16)
    19  = @(510, 12, 484, 14)
    20  = @(505, 8, 457, 40)
    21  = @(501, 22, 451, 28)
    22  = @(485, 3, 447, 35)
    23  = @(483, 42, 403, 38)
    24  = @(480, 72, 459, 21)
    25  = @(477, 4, 467, 6)
    26  = @(438, 8, 394, 21)
    27  = @(430, 8, 390, 32)
    28  = @(427, 72, 413, 13)
    37  = @(257, 5, 248, 4)
    38  = @(245, 72, 242, 3)
    39  = @(234, 72, 228, 6)
    42  = @(189, 25, 186, 3)
    43  = @(179, 14, 152, 13)
    44  = @(158, 8, 156, 21)
    46  = @(103, 25, 99, 4)
    47  = @(100, 11, 82, 7)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - Casos totales
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - Casos activos
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - Recuperados
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - Muertes
}
